$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.884.93'
$ws.Range("E2").Value = '  +0.95%  '

$ws.Range("D3").Value = '1.887.98'
$ws.Range("E3").Value = '  +0.50%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '335.29'
$ws.Range("E5").Value = '  +1.60%  '

$ws.Range("E6").Value = '  +1.59%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4681'
$ws.Range("E7").Value = '  -0.69%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3918'
$ws.Range("E8").Value = '  -1.79%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '47.71'
$ws.Range("E9").Value = '  +1.16%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08047'
$ws.Range("E10").Value = '  -0.27%  '

$ws.Range("E11").Value = '  -0.87%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '21.76'
$ws.Range("E12").Value = '  -0.58%  '

$ws.Range("D13").Value = '1.883.22'
$ws.Range("E13").Value = '  +0.36%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.953'
$ws.Range("E14").Value = '  -0.20%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.096'
$ws.Range("E15").Value = '  -1.65%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.019'
$ws.Range("E16").Value = '  +1.66%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.06784'
$ws.Range("E17").Value = '  +3.20%  '

$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001049'
$ws.Range("E18").Value = '  +0.61%  '

$ws.Range("B19").Value = 'Litecoin'
$ws.Range("C19").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '87.14'
$ws.Range("E19").Value = '  +0.14%  '

$ws.Range("E20").Value = '  -1.03%  '

$ws.Range("E21").Value = '  +1.62%  '

$ws.Range("D22").Value = '27.887.33'
$ws.Range("E22").Value = '  +0.89%  '

$ws.Range("E23").Value = '  -0.51%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.99'
$ws.Range("E24").Value = '  -0.28%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.341'
$ws.Range("E25").Value = '  +1.73%  '

$ws.Range("D26").Value = '2.110.93'
$ws.Range("E26").Value = '  +1.49%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '159.46'
$ws.Range("E27").Value = '  +3.41%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.03'
$ws.Range("E28").Value = '  -1.62%  '

$ws.Range("E29").Value = '  -1.47%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.466'
$ws.Range("E30").Value = '  -1.71%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '121.94'
$ws.Range("E31").Value = '  -0.56%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9657'
$ws.Range("E32").Value = '  +0.78%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09483'
$ws.Range("E33").Value = '  -0.21%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.645'
$ws.Range("E34").Value = '  +1.20%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.402'
$ws.Range("E35").Value = '  -5.14%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.356'
$ws.Range("E36").Value = '  +0.63%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06106'
$ws.Range("E37").Value = '  -0.15%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02249'
$ws.Range("E38").Value = '  -0.42%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.215'
$ws.Range("E39").Value = '  -0.51%  '

$ws.Range("E40").Value = '  -3.02%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5965'
$ws.Range("E41").Value = '  -0.76%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1879'
$ws.Range("E42").Value = '  -1.12%  '

$ws.Range("E43").Value = '  -1.17%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.267'
$ws.Range("E44").Value = '  +1.60%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5674'
$ws.Range("E45").Value = '  -0.56%  '

$ws.Range("E46").Value = '  -0.08%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.400'
$ws.Range("E47").Value = '  -0.33%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.929'
$ws.Range("E48").Value = '  -0.77%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06920'
$ws.Range("E49").Value = '  +1.41%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '113.67'
$ws.Range("E50").Value = '  +3.26%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.068'
$ws.Range("E51").Value = '  +0.19%  '
